$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Relocate the signature-block formatting (currently on rows
#    24-25) onto its new home (rows 29-30) BEFORE row 24 gets reused
#    as a data row and row 25 is cleared away.
# ------------------------------------------------------------------
$ws.Range("B24:C24").Copy()
$ws.Range("B29:C29").PasteSpecial(-4122)
$ws.Range("H24:J24").Copy()
$ws.Range("H29:J29").PasteSpecial(-4122)

$ws.Range("B25:C25").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)
$ws.Range("H25:J25").Copy()
$ws.Range("H30:J30").PasteSpecial(-4122)

$ws.Range("B29:C29").Merge()
$ws.Range("H29:J29").Merge()
$ws.Range("B30:C30").Merge()
$ws.Range("H30:J30").Merge()

# Old footer merges/content no longer apply at rows 24/25.
$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B25:C25").UnMerge()
$ws.Range("H25:J25").UnMerge()
$ws.Range("B25:J25").Clear()

# ------------------------------------------------------------------
# 2) Copy the "last table row" formatting (borders etc, currently on
#    row 19) down onto the new last row (24) before row 19 becomes a
#    normal interior row. Also stamp the normal interior-row look
#    (taken from row 18) onto the newly-added rows 19-23.
# ------------------------------------------------------------------
$ws.Range("H24:J24").ClearContents()
$ws.Range("B19:J19").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B22:J22").PasteSpecial(-4122)
$ws.Range("B23:J23").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Rewrite the worker/debt detail rows 16-24 with the updated data.
# ------------------------------------------------------------------
$rows = @(
    @{ Row=16; Tipo="CC"; Doc="73573860";   Nombre="NORBEY VILLADA VILLADA PAJARO";      Periodo="1810"; Valor=31249; Salario=781242  },
    @{ Row=17; Tipo="CC"; Doc="73162788";   Nombre="MAYER ENRIQUE CANTILLO ALTAMIRANDA"; Periodo="1806"; Valor=67776; Salario=1694400 },
    @{ Row=18; Tipo="CC"; Doc="9095273";    Nombre="JORGE ELIECER ROQUE JIMENEZ";        Periodo="1806"; Valor=50832; Salario=1270800 },
    @{ Row=19; Tipo="CC"; Doc="73146705";   Nombre="BLAS ALBERTO TORRES MAZA";           Periodo="1806"; Valor=50832; Salario=1270800 },
    @{ Row=20; Tipo="CC"; Doc="73156798";   Nombre="GILBERTO JARABA NUÑEZ";              Periodo="1810"; Valor=31249; Salario=781242  },
    @{ Row=21; Tipo="CC"; Doc="73192056";   Nombre="JESUS MARIA CANTILLO ALTAMIRANDA";   Periodo="1812"; Valor=31249; Salario=781242  },
    @{ Row=22; Tipo="CC"; Doc="73192056";   Nombre="JESUS MARIA CANTILLO ALTAMIRANDA";   Periodo="1806"; Valor=67776; Salario=781242  },
    @{ Row=23; Tipo="CC"; Doc="1137195709"; Nombre="JOINER RUIZ BARON";                  Periodo="1810"; Valor=31249; Salario=781242  },
    @{ Row=24; Tipo="CC"; Doc="1002202212"; Nombre="ANTONIO JOSE ZABALETA BARBOZA";      Periodo="1810"; Valor=31249; Salario=781242  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Tipo
    $ws.Cells.Item($row, 3).Value = $r.Doc
    $ws.Cells.Item($row, 4).Value = $r.Nombre
    $ws.Cells.Item($row, 5).Value = $r.Periodo
    $ws.Cells.Item($row, 6).Value = $r.Valor
    $ws.Cells.Item($row, 7).Value = $r.Salario
}

# ------------------------------------------------------------------
# 4) Update the summary fields above the table.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 393461
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 5) Write the signature/footer block text into its new rows 29-30.
# ------------------------------------------------------------------
$ws.Cells.Item(29, 2).Value = "___________________________________"
$ws.Cells.Item(29, 8).Value = "___________________________________"
$ws.Cells.Item(30, 2).Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(30, 8).Value = "FIRMA DEL REPRESENTANTE LEGAL"
